# "Generate Report for Handback"
# Refresh the handoff/handback timestamps for the 4cec6f04-... file row
# (row 2) on each localized-language sheet, then roll the newest of those
# timestamps up into the Overview sheet's "Latest HO Xliff Generate Date"
# column for that same file.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-08-14 01:05:15"
$zhcn.Range("K2").Value = "2016-08-14 01:05:44"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-08-14 01:05:23"
$dede.Range("K2").Value = "2016-08-14 01:05:53"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-08-14 01:05:23"
